$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2 through 425 all get updated from serial 45182 (2023-09-13)
# to serial 45184 (2023-09-15), preserving existing date number formatting/style.
$ws.Range("C2:C425").Value = 45184
